$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Structural edits: insert the URL formula column (before old column C) ---
$ws.Columns.Item(3).EntireColumn.Insert()

# --- Insert a new row for the seats/<int:pk> PUT entry (pushes old DELETE row down) ---
$ws.Rows.Item(18).EntireRow.Insert()

# --- Cell content for every row (A,B,D..I); column C (URL) formulas are set afterwards ---
$ws.Range("A1").Value = 'View-name'
$ws.Range("B1").Value = 'Method'
$ws.Range("D1").Value = 'Resouce'
$ws.Range("H1").Value = 'Permission'
$ws.Range("I1").Value = '비고'

$ws.Range("A2").Value = 'student-list'
$ws.Range("B2").Value = 'GET'
$ws.Range("D2").Value = 'students'
$ws.Range("H2").Value = 'IsAdminUser'
$ws.Range("I2").Value = '(create, destroy는 자동)'

$ws.Range("A3").Value = 'student-detail'
$ws.Range("B3").Value = 'GET'
$ws.Range("D3").Value = 'students'
$ws.Range("E3").Value = '<int:pk>'
$ws.Range("H3").Value = 'IsOwner | IsAdminUser'

$ws.Range("A4").Value = 'student-purchase'
$ws.Range("B4").Value = 'GET'
$ws.Range("D4").Value = 'students'
$ws.Range("E4").Value = '<int:pk>'
$ws.Range("F4").Value = 'purchases'
$ws.Range("H4").Value = 'IsOwner | IsAdminUser'

$ws.Range("A5").Value = 'student-rents'
$ws.Range("B5").Value = 'GET'
$ws.Range("D5").Value = 'students'
$ws.Range("E5").Value = '<int:pk>'
$ws.Range("F5").Value = 'rents'
$ws.Range("G5").Value = '?year=2020&month=4'
$ws.Range("H5").Value = 'IsOwner | IsAdminUser'

$ws.Range("A6").Value = '학생의 현재 이용권 저장가능'
$ws.Range("B6").Value = 'GET'
$ws.Range("D6").Value = 'students'
$ws.Range("E6").Value = '<int:pk>'
$ws.Range("F6").Value = 'ticket-storable'
$ws.Range("H6").Value = 'IsOwner | IsAdminUser'

$ws.Range("A7").Value = 'student-update'
$ws.Range("B7").Value = 'PUT'
$ws.Range("D7").Value = 'students'
$ws.Range("E7").Value = '<int:pk>'

$ws.Range("A9").Value = 'ticket-list'
$ws.Range("B9").Value = 'GET'
$ws.Range("D9").Value = 'tickets'
$ws.Range("H9").Value = 'AllowAny'

$ws.Range("A10").Value = 'ticket-detail'
$ws.Range("B10").Value = 'GET'
$ws.Range("D10").Value = 'tickets'
$ws.Range("E10").Value = '<int:pk>'
$ws.Range("H10").Value = 'AllowAny'

$ws.Range("A11").Value = 'ticket-create'
$ws.Range("B11").Value = 'POST'
$ws.Range("D11").Value = 'tickets'
$ws.Range("H11").Value = 'IsAdminUser'

$ws.Range("A12").Value = 'ticket-update'
$ws.Range("B12").Value = 'PUT'
$ws.Range("D12").Value = 'tickets'
$ws.Range("E12").Value = '<int:pk>'
$ws.Range("H12").Value = 'IsAdminUser'
$ws.Range("I12").Value = '(구매 없는 것만 허용)'

$ws.Range("A13").Value = 'ticket-destroy'
$ws.Range("B13").Value = 'DELETE'
$ws.Range("D13").Value = 'tickets'
$ws.Range("E13").Value = '<int:pk>'
$ws.Range("H13").Value = 'IsAdminUser'
$ws.Range("I13").Value = '(구매 없는 것만 허용)'

$ws.Range("A15").Value = 'seat-list'
$ws.Range("B15").Value = 'GET'
$ws.Range("D15").Value = 'seats'
$ws.Range("H15").Value = 'AllowAny'
$ws.Range("I15").Value = '(seat 필드 추가 시 seat-detail, seat-update 추가)'

$ws.Range("A16").Value = '대여 가능 자리 목록'
$ws.Range("B16").Value = 'GET'
$ws.Range("D16").Value = 'seats'
$ws.Range("E16").Value = '?'
$ws.Range("H16").Value = 'AllowAny'
$ws.Range("I16").Value = '(미완)'

$ws.Range("A17").Value = 'seat-create'
$ws.Range("B17").Value = 'POST'
$ws.Range("D17").Value = 'seats'
$ws.Range("H17").Value = 'IsAdminUser'

$ws.Range("A18").Value = 'seat-update'
$ws.Range("B18").Value = 'PUT'
$ws.Range("D18").Value = 'seats'
$ws.Range("E18").Value = '<int:pk>'
$ws.Range("H18").Value = 'IsAdminUser'
$ws.Range("I18").Value = '(대여 없는 것만 허용)'

$ws.Range("A19").Value = 'seat-update'
$ws.Range("B19").Value = 'DELETE'
$ws.Range("D19").Value = 'seats'
$ws.Range("E19").Value = '<int:pk>'
$ws.Range("H19").Value = 'IsAdminUser'
$ws.Range("I19").Value = '(대여 없는 것만 허용)'

$ws.Range("A21").Value = 'purchase-list'
$ws.Range("B21").Value = 'GET'
$ws.Range("D21").Value = 'purchases'
$ws.Range("H21").Value = 'IsAdminUser'

$ws.Range("A22").Value = '월 매출'
$ws.Range("B22").Value = 'GET'
$ws.Range("D22").Value = 'purchases'
$ws.Range("E22").Value = 'price'
$ws.Range("F22").Value = '?year=2020'
$ws.Range("H22").Value = 'IsAdminUser'

$ws.Range("A23").Value = 'purchase-create'
$ws.Range("B23").Value = 'POST'
$ws.Range("D23").Value = 'purchases'
$ws.Range("H23").Value = 'IsAuthenticated'

$ws.Range("A25").Value = 'rent-list'
$ws.Range("B25").Value = 'GET'
$ws.Range("D25").Value = 'rents'
$ws.Range("H25").Value = 'IsAdminUser'

$ws.Range("A26").Value = 'rent-create (start)'
$ws.Range("B26").Value = 'POST'
$ws.Range("D26").Value = 'rents'
$ws.Range("H26").Value = 'IsAuthenticated'

$ws.Range("A27").Value = 'rent-update (end)'
$ws.Range("B27").Value = 'PUT'
$ws.Range("D27").Value = 'rents'
$ws.Range("H27").Value = 'IsAuthenticated'

$ws.Range("A29").Value = 'user-list'
$ws.Range("B29").Value = 'GET'
$ws.Range("D29").Value = 'users'
$ws.Range("H29").Value = 'AllowAny'

$ws.Range("A30").Value = 'user-detail'
$ws.Range("B30").Value = 'GET'
$ws.Range("D30").Value = 'users'
$ws.Range("E30").Value = '<int:pk>'
$ws.Range("H30").Value = 'IsOwner | IsAdminUser'

$ws.Range("A31").Value = 'user-create (signup)'
$ws.Range("B31").Value = 'POST'
$ws.Range("D31").Value = 'users'
$ws.Range("H31").Value = 'AllowAny'

$ws.Range("A32").Value = 'user-update (email, student)'
$ws.Range("B32").Value = 'PUT'
$ws.Range("D32").Value = 'users'
$ws.Range("E32").Value = '<int:pk>'
$ws.Range("H32").Value = 'IsOwner'

$ws.Range("A33").Value = 'user password change'
$ws.Range("B33").Value = 'PUT'
$ws.Range("D33").Value = 'users'
$ws.Range("E33").Value = '<int:pk>'
$ws.Range("F33").Value = 'password'
$ws.Range("H33").Value = 'IsOwner'

$ws.Range("A34").Value = 'user-destroy'
$ws.Range("B34").Value = 'DELETE'
$ws.Range("D34").Value = 'users'
$ws.Range("E34").Value = '<int:pk>'
$ws.Range("H34").Value = 'IsOwner'
$ws.Range("I34").Value = '(구매 없는 것만 허용?)'

# --- Column C: URL formula (row 1 is a literal header string, not a formula) ---
$ws.Range("C1").Value = 'URL'
$ws.Range("C2").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D2),"/"&D2,""),IF(ISTEXT(E2),"/"&E2,""),IF(ISTEXT(F2),"/"&F2,""),IF(ISTEXT(G2),"/"&G2,""))'
$ws.Range("C3").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D3),"/"&D3,""),IF(ISTEXT(E3),"/"&E3,""),IF(ISTEXT(F3),"/"&F3,""),IF(ISTEXT(G3),"/"&G3,""))'
$ws.Range("C4").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D4),"/"&D4,""),IF(ISTEXT(E4),"/"&E4,""),IF(ISTEXT(F4),"/"&F4,""),IF(ISTEXT(G4),"/"&G4,""))'
$ws.Range("C5").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D5),"/"&D5,""),IF(ISTEXT(E5),"/"&E5,""),IF(ISTEXT(F5),"/"&F5,""),IF(ISTEXT(G5),"/"&G5,""))'
$ws.Range("C6").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D6),"/"&D6,""),IF(ISTEXT(E6),"/"&E6,""),IF(ISTEXT(F6),"/"&F6,""),IF(ISTEXT(G6),"/"&G6,""))'
$ws.Range("C7").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D7),"/"&D7,""),IF(ISTEXT(E7),"/"&E7,""),IF(ISTEXT(F7),"/"&F7,""),IF(ISTEXT(G7),"/"&G7,""))'
$ws.Range("C8").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D8),"/"&D8,""),IF(ISTEXT(E8),"/"&E8,""),IF(ISTEXT(F8),"/"&F8,""),IF(ISTEXT(G8),"/"&G8,""))'
$ws.Range("C9").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D9),"/"&D9,""),IF(ISTEXT(E9),"/"&E9,""),IF(ISTEXT(F9),"/"&F9,""),IF(ISTEXT(G9),"/"&G9,""))'
$ws.Range("C10").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D10),"/"&D10,""),IF(ISTEXT(E10),"/"&E10,""),IF(ISTEXT(F10),"/"&F10,""),IF(ISTEXT(G10),"/"&G10,""))'
$ws.Range("C11").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D11),"/"&D11,""),IF(ISTEXT(E11),"/"&E11,""),IF(ISTEXT(F11),"/"&F11,""),IF(ISTEXT(G11),"/"&G11,""))'
$ws.Range("C12").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D12),"/"&D12,""),IF(ISTEXT(E12),"/"&E12,""),IF(ISTEXT(F12),"/"&F12,""),IF(ISTEXT(G12),"/"&G12,""))'
$ws.Range("C13").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D13),"/"&D13,""),IF(ISTEXT(E13),"/"&E13,""),IF(ISTEXT(F13),"/"&F13,""),IF(ISTEXT(G13),"/"&G13,""))'
$ws.Range("C14").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D14),"/"&D14,""),IF(ISTEXT(E14),"/"&E14,""),IF(ISTEXT(F14),"/"&F14,""),IF(ISTEXT(G14),"/"&G14,""))'
$ws.Range("C15").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D15),"/"&D15,""),IF(ISTEXT(E15),"/"&E15,""),IF(ISTEXT(F15),"/"&F15,""),IF(ISTEXT(G15),"/"&G15,""))'
$ws.Range("C16").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D16),"/"&D16,""),IF(ISTEXT(E16),"/"&E16,""),IF(ISTEXT(F16),"/"&F16,""),IF(ISTEXT(G16),"/"&G16,""))'
$ws.Range("C17").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D17),"/"&D17,""),IF(ISTEXT(E17),"/"&E17,""),IF(ISTEXT(F17),"/"&F17,""),IF(ISTEXT(G17),"/"&G17,""))'
$ws.Range("C18").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D18),"/"&D18,""),IF(ISTEXT(E18),"/"&E18,""),IF(ISTEXT(F18),"/"&F18,""),IF(ISTEXT(G18),"/"&G18,""))'
$ws.Range("C19").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D19),"/"&D19,""),IF(ISTEXT(E19),"/"&E19,""),IF(ISTEXT(F19),"/"&F19,""),IF(ISTEXT(G19),"/"&G19,""))'
$ws.Range("C20").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D20),"/"&D20,""),IF(ISTEXT(E20),"/"&E20,""),IF(ISTEXT(F20),"/"&F20,""),IF(ISTEXT(G20),"/"&G20,""))'
$ws.Range("C21").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D21),"/"&D21,""),IF(ISTEXT(E21),"/"&E21,""),IF(ISTEXT(F21),"/"&F21,""),IF(ISTEXT(G21),"/"&G21,""))'
$ws.Range("C22").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D22),"/"&D22,""),IF(ISTEXT(E22),"/"&E22,""),IF(ISTEXT(F22),"/"&F22,""),IF(ISTEXT(G22),"/"&G22,""))'
$ws.Range("C23").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D23),"/"&D23,""),IF(ISTEXT(E23),"/"&E23,""),IF(ISTEXT(F23),"/"&F23,""),IF(ISTEXT(G23),"/"&G23,""))'
$ws.Range("C24").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D24),"/"&D24,""),IF(ISTEXT(E24),"/"&E24,""),IF(ISTEXT(F24),"/"&F24,""),IF(ISTEXT(G24),"/"&G24,""))'
$ws.Range("C25").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D25),"/"&D25,""),IF(ISTEXT(E25),"/"&E25,""),IF(ISTEXT(F25),"/"&F25,""),IF(ISTEXT(G25),"/"&G25,""))'
$ws.Range("C26").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D26),"/"&D26,""),IF(ISTEXT(E26),"/"&E26,""),IF(ISTEXT(F26),"/"&F26,""),IF(ISTEXT(G26),"/"&G26,""))'
$ws.Range("C27").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D27),"/"&D27,""),IF(ISTEXT(E27),"/"&E27,""),IF(ISTEXT(F27),"/"&F27,""),IF(ISTEXT(G27),"/"&G27,""))'
$ws.Range("C28").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D28),"/"&D28,""),IF(ISTEXT(E28),"/"&E28,""),IF(ISTEXT(F28),"/"&F28,""),IF(ISTEXT(G28),"/"&G28,""))'
$ws.Range("C29").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D29),"/"&D29,""),IF(ISTEXT(E29),"/"&E29,""),IF(ISTEXT(F29),"/"&F29,""),IF(ISTEXT(G29),"/"&G29,""))'
$ws.Range("C30").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D30),"/"&D30,""),IF(ISTEXT(E30),"/"&E30,""),IF(ISTEXT(F30),"/"&F30,""),IF(ISTEXT(G30),"/"&G30,""))'
$ws.Range("C31").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D31),"/"&D31,""),IF(ISTEXT(E31),"/"&E31,""),IF(ISTEXT(F31),"/"&F31,""),IF(ISTEXT(G31),"/"&G31,""))'
$ws.Range("C32").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D32),"/"&D32,""),IF(ISTEXT(E32),"/"&E32,""),IF(ISTEXT(F32),"/"&F32,""),IF(ISTEXT(G32),"/"&G32,""))'
$ws.Range("C33").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D33),"/"&D33,""),IF(ISTEXT(E33),"/"&E33,""),IF(ISTEXT(F33),"/"&F33,""),IF(ISTEXT(G33),"/"&G33,""))'
$ws.Range("C34").Formula = '=_xlfn.CONCAT(IF(ISTEXT(D34),"/"&D34,""),IF(ISTEXT(E34),"/"&E34,""),IF(ISTEXT(F34),"/"&F34,""),IF(ISTEXT(G34),"/"&G34,""))'

# --- Styling: row 7 (student-update) gray/disabled look, matches legacy 'destroy is automatic' rows ---
$ws.Range("A7").Interior.Color = 65535
$ws.Range("A7").Font.Color = 8421504
$ws.Range("B7").Font.Color = 8421504
$ws.Range("D7").Font.Color = 8421504
$ws.Range("E7").Font.Color = 8421504

# --- Column widths / visibility ---
$ws.Columns.Item(3).ColumnWidth = 16.5
$ws.Range("D1:F1").EntireColumn.Hidden = $true
$ws.Columns.Item(7).ColumnWidth = 20.8984375
$ws.Columns.Item(7).EntireColumn.Hidden = $true
$ws.Columns.Item(8).ColumnWidth = 11.8984375

# --- View state: selection + scroll position ---
$ws.Application.GoTo($ws.Range("A25"))
$ws.Range("C2:C7").Select()
